# FsStats.xlsx edit: "Aggiunto boxplot e medie"
#  - Rename "Lirlanda" -> "Lirlanda FC" (and the combined "Lirlanda, SS Egli
#    Tare 2016" -> "Lirlanda FC, SS Egli Tare 2016") on the "Albo" sheet.
#  - Move/relabel the little underline marker cell on "Albo" from D14 to F14
#    (new underlined font), widen column E a bit.
#  - Add a new "Riassunto" summary sheet at the end of the workbook.
#  - Restore cursor/selection positions + the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Albo" sheet: rename "Lirlanda" -> "Lirlanda FC" everywhere it shows
#    up (the club added "FC" to its name).
# ---------------------------------------------------------------------
$albo = $wb.Worksheets.Item("Albo")

$albo.Range("B2").Value = "Lirlanda FC"
$albo.Range("D4").Value = "Lirlanda FC"
$albo.Range("C8").Value = "Lirlanda FC"
# NB: E9 ("Lirlanda FC, SS Egli Tare 2016") is set further down, after the
# new sheet's strings are created, to match the shared-string table order
# of the target workbook (new strings are appended in first-use order).

# Widen column E (Retrocessioni) now that the longer name needs more room.
$albo.Columns.Item(5).ColumnWidth = 25.859375

# Move the little formatting-only marker cell from D14 to F14, giving it
# an underlined font.
$albo.Cells.Item(14, 4).Font.Underline = 2
$albo.Cells.Item(14, 4).Cut($albo.Cells.Item(14, 6)) | Out-Null
$albo.Cells.Item(14, 4).Clear() | Out-Null

$albo.Range("E16").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) "16-17" sheet: just the remembered cursor position changed.
# ---------------------------------------------------------------------
$s1617 = $wb.Worksheets.Item("16-17")
$s1617.Range("U29").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) Add the new "Riassunto" (summary) sheet after the last sheet ("Albo").
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$riassunto = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$riassunto.Name = "Riassunto"

$riassunto.Range("B1").Value = "16-17"
$riassunto.Range("C1").Value = "17-18"
$riassunto.Range("D1").Value = "18-19"
$riassunto.Range("E1").Value = "19-20"
$riassunto.Range("F1").Value = "20-21"
$riassunto.Range("G1").Value = "21-22"
$riassunto.Range("H1").Value = "22-23"
$riassunto.Range("I1").Value = "23-24"

$riassunto.Range("A2").Value = "Calendario"
$riassunto.Range("B2").Value = "Si"
$riassunto.Range("C2").Value = "Si"
$riassunto.Range("D2").Value = "*"
$riassunto.Range("E2").Value = "-"
$riassunto.Range("F2").Value = "Si"
$riassunto.Range("G2").Value = "Si"
$riassunto.Range("H2").Value = "Si"
$riassunto.Range("I2").Value = "Si"

$riassunto.Range("A3").Value = "Rose"
$riassunto.Range("B3").Value = "Si"
$riassunto.Range("C3").Value = "Si"
$riassunto.Range("D3").Value = "Si"
$riassunto.Range("E3").Value = "Si"
$riassunto.Range("F3").Value = "Si"
$riassunto.Range("G3").Value = "-"
$riassunto.Range("H3").Value = "-"
$riassunto.Range("I3").Value = "Si"

$riassunto.Range("A4").Value = "Dettagli"
$riassunto.Range("B4").Value = "Si"
$riassunto.Range("C4").Value = "Si"
$riassunto.Range("D4").Value = "-"
$riassunto.Range("E4").Value = "Si"
$riassunto.Range("F4").Value = "-"
$riassunto.Range("G4").Value = "-"
$riassunto.Range("H4").Value = "-"
$riassunto.Range("I4").Value = "-"

$riassunto.Range("A6").Value = "*"
$riassunto.Range("B6").Value = "Abbiamo le formazioni,  è recuperabile con un po' di lavoro"

$riassunto.Range("I9").Select() | Out-Null

# Now that the "Riassunto" strings exist, set the last renamed cell so the
# shared-string table ends up in the same append order as the target file.
$albo.Range("E9").Value = "Lirlanda FC, SS Egli Tare 2016"

# ---------------------------------------------------------------------
# 4) Restore "Albo" as the active/selected tab (adding the sheet above
#    made the new sheet active).
# ---------------------------------------------------------------------
$albo.Activate()
